$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '25.954.75'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = "'" + '1.641.71'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").Value = "'" + '215.16'
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").Value = "'" + '0.5065'
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("D7").Value = "'" + '1.003'
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = "'" + '0.2557'
$ws.Range("E8").Value = '  -0.55%  '
$ws.Range("D9").Value = "'" + '0.06373'
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("D10").Value = "'" + '19.46'
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("D11").Value = "'" + '0.07757'
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = "'" + '4.278'
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = "'" + '1.646.01'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("D14").Value = "'" + '0.5449'
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("D15").Value = "'" + '0.0' + [char]0x2085 + '7822'
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").Value = "'" + '64.25'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = "'" + '25.993.64'
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").Value = "'" + '197.57'
$ws.Range("E19").Value = '  -2.14%  '
$ws.Range("D20").Value = "'" + '4.439'
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("D21").Value = "'" + '9.938'
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").Value = "'" + '6.047'
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("D24").Value = "'" + '1.896'
$ws.Range("E24").Value = '  +1.72%  '
$ws.Range("D25").Value = "'" + '140.85'
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E26").Value = '  +2.89%  '
$ws.Range("D27").Value = "'" + '6.884'
$ws.Range("E27").Value = '  +1.63%  '
$ws.Range("D28").Value = "'" + '15.70'
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("D29").Value = "'" + '1.237'
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("D30").Value = "'" + '0.04958'
$ws.Range("E30").Value = '  +0.68%  '
$ws.Range("D31").Value = "'" + '3.258'
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("D32").Value = "'" + '3.183'
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("D33").Value = "'" + '1.538'
$ws.Range("E33").Value = '  -0.47%  '
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").Value = "'" + '0.8951'
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("D37").Value = "'" + '1.133.69'
$ws.Range("E37").Value = '  -1.03%  '
$ws.Range("D38").Value = "'" + '0.5442'
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("D39").Value = "'" + '0.01560'
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("B40").Value = 'mCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D40").Value = "'" + '2.556'
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = "'" + '1.003'
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").Value = "'" + '0.8201'
$ws.Range("E42").Value = '  +1.88%  '
$ws.Range("D43").Value = "'" + '5.582'
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("D44").Value = "'" + '0.0' + [char]0x2088 + '127'
$ws.Range("E44").Value = '  +7.50%  '
$ws.Range("D45").Value = "'" + '99.51'
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").Value = "'" + '1.777.59'
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").Value = "'" + '0.4536'
$ws.Range("E47").Value = '  +0.64%  '
$ws.Range("D48").Value = "'" + '1.002'
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").Value = "'" + '54.72'
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").Value = "'" + '0.05071'
$ws.Range("E50").Value = '  +0.51%  '
$ws.Range("D51").Value = "'" + '1.005'
$ws.Range("E51").Value = '  +0.35%  '
